$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 15994
$ws1.Range("F9").Value = 15511
$ws1.Range("F11").Value = 9138
$ws1.Range("F12").Value = 407
$ws1.Range("F14").Value = 1022
$ws1.Range("F15").Value = 112
$ws1.Range("F18").Value = 212
$ws1.Range("F21").Value = 576
$ws1.Range("F26").Value = 10
$ws1.Range("F29").Value = 505
$ws1.Range("F30").Value = 36
$ws1.Range("F34").Value = 57
$ws1.Range("F36").Value = 335
$ws1.Range("F37").Value = 463
$ws1.Range("F39").Value = 5618
$ws1.Range("F40").Value = 5234

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 15994
$ws4.Range("F9").Value = 15511
$ws4.Range("F11").Value = 9138
$ws4.Range("F12").Value = 407
$ws4.Range("F14").Value = 1022
$ws4.Range("F15").Value = 112
$ws4.Range("F18").Value = 212
$ws4.Range("F21").Value = 576
$ws4.Range("F26").Value = 10
$ws4.Range("F29").Value = 505
$ws4.Range("F30").Value = 36
$ws4.Range("F36").Value = 57
$ws4.Range("F38").Value = 335
$ws4.Range("F39").Value = 463
$ws4.Range("F41").Value = 5618
$ws4.Range("F43").Value = 5234
